$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.728.66"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "2.214.50"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'229.45"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "'0.634"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").Value = "'64.21"
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.405"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "2.542.30"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "'15.88"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "'22.36"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "2.215.24"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "40.567.15"
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("D19").Value = "'73.87"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("E20").Value = "  +6.08%  "
$ws.Range("D21").Value = "'6.12"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'250.17"
$ws.Range("E22").Value = "  +7.64%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").Value = "'9.74"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").Value = "'173.04"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Value = "'20.41"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").Value = "'4.68"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").Value = "'4.77"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").Value = "'7.13"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'0.0630"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").Value = "'3.83"
$ws.Range("E37").Value = "  +6.24%  "
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'4.79"
$ws.Range("E40").Value = "  +9.15%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "'8.56"
$ws.Range("E42").Value = "  +8.64%  "
$ws.Range("D43").Value = "'101.45"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("D45").Value = "1.521.32"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").Value = "'17.31"
$ws.Range("E46").Value = "  -3.62%  "
$ws.Range("D47").Value = "'0.0938"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.11"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("B49").Value = "TerraClassic"
$ws.Range("C49").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D49").Value = "'0.000207"
$ws.Range("E49").Value = "  +39.92%  "
$ws.Range("D50").Value = "'2.82"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "'50.99"
$ws.Range("E51").Value = "  +9.47%  "
